$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the batch-size related values / SQL strings (null->false fix changed
# the expected boundary values from 400/6000/5600 to 4000/20000/16000).
$ws.Range("I2").Value = "20000"
$ws.Range("L2").Value = "20000"
$ws.Range("T2").Value = "20000"

$ws.Range("M2").Value = "select count(*) from `$schema26 where id>4000 and id<=20000"

$ws.Range("N2").Value = "16000"
$ws.Range("P2").Value = "16000"
$ws.Range("R2").Value = "16000"

$ws.Range("O2").Value = "update `$schema26 set name='BJ' where id>4000 and id<=20000"

# Refresh the view/selection state to match the authored commit (scrolled
# right so column M is the left-most visible column, with V2 selected).
$ws.Activate()
$ws.Range("V2").Select()
$excel.ActiveWindow.ScrollColumn = 13
